$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Remove the old "Tests for Factory" block (rows 17:19) - it is being
# replaced by the new "Marshal by reference" block (rows 17:26).
# ---------------------------------------------------------------------
$ws.Range("A17:M19").ClearContents()

# ---------------------------------------------------------------------
# Column widths: drop the old bestFit column A width and instead size
# columns A and D explicitly.
# ---------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 10
$ws.Columns.Item(4).ColumnWidth = 18

# ---------------------------------------------------------------------
# The overall AND() check at the top now covers through row 26.
# ---------------------------------------------------------------------
$ws.Range("B1").Formula = "=AND(A4:A26)"

# ---------------------------------------------------------------------
# Row 17: section headers (bold)
# ---------------------------------------------------------------------
$ws.Range("A17").Value = "Marshal by reference"
$ws.Range("A17").Font.Bold = $true

$ws.Range("D17").Value = "One call returning IEnumerable<T>"
$ws.Range("D17").Font.Bold = $true

$ws.Range("G17").Value = "Multiple calls returning T"
$ws.Range("G17").Font.Bold = $true

# ---------------------------------------------------------------------
# First block: rows 18-20
# ---------------------------------------------------------------------
$ws.Range("D18:D20").FormulaArray = "=_xll.dnaFactoryMultiple(E18:E20,F18:F20)"

$ws.Range("E18").Value = "One"
$ws.Range("F18").Value = 1
$ws.Range("G18").Formula = "=_xll.dnaFactorySingle(E18,F18)"

$ws.Range("E19").Value = "Two"
$ws.Range("F19").Value = 2
$ws.Range("G19").Formula = "=_xll.dnaFactorySingle(E19,F19)"
$ws.Range("J19").Formula = "=_xll.dnaFactoryCompound(G19,G18)"

$ws.Range("E20").Value = "One"
$ws.Range("F20").Value = 3
$ws.Range("G20").Formula = "=_xll.dnaFactorySingle(E20,F20)"
$ws.Range("J20").Formula = "=_xll.dnaFactoryCompound(J19,G20)"

# Row 21: marshal-by-ref summary for the first block
$ws.Range("A21").Formula = "=SUMSQ(B21:C21)<0.00000000000001"
$ws.Range("B21").Formula = "=D21-G21"
$ws.Range("C21").Formula = "=G21-J21"
$ws.Range("D21").Formula = "=_xll.dnaMarshalByRef(D18:D20)"
$ws.Range("G21").Formula = "=_xll.dnaMarshalByRef(G18:G20)"
$ws.Range("J21").Formula = "=_xll.dnaMarshalByRef(J20)"

# ---------------------------------------------------------------------
# Second block: rows 23-25
# ---------------------------------------------------------------------
$ws.Range("D23:D25").FormulaArray = "=_xll.dnaFactoryMultiple(E23:E25,F23:F25)"

$ws.Range("E23").Value = "Two"
$ws.Range("F23").Value = 4
$ws.Range("G23").Formula = "=_xll.dnaFactorySingle(E23,F23)"

$ws.Range("E24").Value = "One"
$ws.Range("F24").Value = 5
$ws.Range("G24").Formula = "=_xll.dnaFactorySingle(E24,F24)"
$ws.Range("J24").Formula = "=_xll.dnaFactoryCompound(G24,G23)"

$ws.Range("E25").Value = "Two"
$ws.Range("F25").Value = 6
$ws.Range("G25").Formula = "=_xll.dnaFactorySingle(E25,F25)"
$ws.Range("J25").Formula = "=_xll.dnaFactoryCompound(J24,G25)"

# Row 26: marshal-by-ref summary for the second block
$ws.Range("A26").Formula = "=SUMSQ(B26:C26)<0.00000000000001"
$ws.Range("B26").Formula = "=D26-G26"
$ws.Range("C26").Formula = "=G26-J26"
$ws.Range("D26").Formula = "=_xll.dnaMarshalByRef(D23:D25)"
$ws.Range("G26").Formula = "=_xll.dnaMarshalByRef(G23:G25)"
$ws.Range("J26").Formula = "=_xll.dnaMarshalByRef(J25)"

# ---------------------------------------------------------------------
# Selection moves to B1 (matches the new authored state of the sheet).
# ---------------------------------------------------------------------
$ws.Range("B1").Select()
